$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 228, shifting all existing rows (228..244) down to (229..245)
$ws.Rows.Item(228).Insert()

# Populate the newly inserted row 228 with the new record
$ws.Cells.Item(228, 1).Value = 7
$ws.Cells.Item(228, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(228, 3).Value = "Ñuble"
$ws.Cells.Item(228, 4).Value = 44783
$ws.Cells.Item(228, 5).Value = 16
$ws.Cells.Item(228, 6).Value = 100112003
$ws.Cells.Item(228, 7).Value = "Ajo"
$ws.Cells.Item(228, 8).Value = "Chino"
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 60
$ws.Cells.Item(228, 11).Value = 24000
$ws.Cells.Item(228, 12).Value = 25000
$ws.Cells.Item(228, 13).Value = 24500
$ws.Cells.Item(228, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(228, 15).Value = "China"
$ws.Cells.Item(228, 16).Value = 2450
$ws.Cells.Item(228, 17).Value = 10
$ws.Cells.Item(228, 18).Value = "Hortaliza"
